$d = $word.ActiveDocument

function Wrap-Xml($inner) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $inner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Work from the bottom of the document upward so paragraph indices for
# not-yet-processed paragraphs stay stable while we insert/remove paragraphs.

# --- Paragraphs 44 & 45: merge "List custom levels in Add-ons menu?" and
#     "List custom levels in <<Bonus Levels>>?" into a single struck-through
#     "List custom levels in <<Custom Levels>>?" paragraph.
$p44 = $d.Paragraphs(44)
$p45 = $d.Paragraphs(45)
$rng = $d.Range($p44.Range.Start, $p45.Range.End)
$inner = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t>List custom levels in “</w:t></w:r>' + `
    '<w:r><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t>Custom Levels</w:t></w:r>' + `
    '<w:r><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t>”?</w:t></w:r>' + `
    '</w:p>'
$rng.InsertXML((Wrap-Xml $inner))

# --- Paragraph 43: "How will users play custom levels outside of the level editor?" -> add strike
$p = $d.Paragraphs(43)
$inner = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t>How will users play custom levels outside of the level editor?</w:t></w:r>' + `
    '</w:p>'
$p.Range.InsertXML((Wrap-Xml $inner))

# --- Paragraph 41: "Yes. You can simply open them in the editor" -> add strike
$p = $d.Paragraphs(41)
$inner = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Yes. You can simply open them in the </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/><w:r><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t>editor</w:t></w:r><w:proofErr w:type="gramEnd"/>' + `
    '</w:p>'
$p.Range.InsertXML((Wrap-Xml $inner))

# --- Paragraph 40: "Should individual levels be accessible without a worldmap?"
#     -> split off "worldmap" with spell-check markers (already struck through)
$p = $d.Paragraphs(40)
$inner = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Should individual levels be accessible without a </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t>worldmap</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t>?</w:t></w:r>' + `
    '</w:p>'
$p.Range.InsertXML((Wrap-Xml $inner))

# --- Paragraph 38: "Do levels need a worldmap to be accessible?" -> add strike,
#     split off "worldmap" with spell-check markers, then insert new "Yes." paragraph after it.
$p = $d.Paragraphs(38)
$inner = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Do </w:t></w:r>' + `
    '<w:r><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t>levels</w:t></w:r>' + `
    '<w:r><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> need a </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t>worldmap</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> to be accessible?</w:t></w:r>' + `
    '</w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t>Yes.</w:t></w:r>' + `
    '</w:p>'
$p.Range.InsertXML((Wrap-Xml $inner))

# --- Paragraph 37 (empty paragraph) -> add strike
$p = $d.Paragraphs(37)
$inner = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>'
$p.Range.InsertXML((Wrap-Xml $inner))

# --- Paragraph 36: "How will worldmaps ensure that the levels within them have
#     valid file paths?" -> add strike, split off "worldmaps" with spell-check
#     markers, then insert new "Use a custom selector UI..." paragraph after it.
$p = $d.Paragraphs(36)
$inner = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">How will </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t>worldmaps</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> ensure that the levels within them have valid file paths?</w:t></w:r>' + `
    '</w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t>Use a custom selector UI that only loads levels in the world.</w:t></w:r>' + `
    '</w:p>'
$p.Range.InsertXML((Wrap-Xml $inner))

# --- Paragraph 30: "Ability to load / open levels from disk" -> add strike
$p = $d.Paragraphs(30)
$inner = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/><w:bCs/><w:strike/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:b/><w:bCs/><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t>Ability to load / open levels</w:t></w:r>' + `
    '<w:r><w:rPr><w:b/><w:bCs/><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> from </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t>disk</w:t></w:r><w:proofErr w:type="gramEnd"/>' + `
    '</w:p>'
$p.Range.InsertXML((Wrap-Xml $inner))

# --- Paragraph 29: "Ability to save levels to disk" -> add strike
$p = $d.Paragraphs(29)
$inner = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/><w:bCs/><w:strike/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:b/><w:bCs/><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t>Ability to save levels</w:t></w:r>' + `
    '<w:r><w:rPr><w:b/><w:bCs/><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> to </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t>disk</w:t></w:r><w:proofErr w:type="gramEnd"/>' + `
    '</w:p>'
$p.Range.InsertXML((Wrap-Xml $inner))

# --- Paragraph 27: "Ability to change the level properties (name, author, gravity, music)" -> add strike
$p = $d.Paragraphs(27)
$inner = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t>Ability to change the level properties (name, author, gravity, music)</w:t></w:r>' + `
    '</w:p>'
$p.Range.InsertXML((Wrap-Xml $inner))

# --- Paragraph 26: "Ability to change the background" -> add strike
$p = $d.Paragraphs(26)
$inner = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Ability to change the </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/><w:r><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t>background</w:t></w:r><w:proofErr w:type="gramEnd"/>' + `
    '</w:p>'
$p.Range.InsertXML((Wrap-Xml $inner))

# --- Paragraph 23: "Ability to add or remove TileMaps" -> wrap "TileMaps" with spell-check markers
$p = $d.Paragraphs(23)
$inner = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Ability to add or remove </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t>TileMaps</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' + `
    '</w:p>'
$p.Range.InsertXML((Wrap-Xml $inner))

# --- Paragraph 22: "Ability to change the active TileMap" -> wrap "TileMap" with spell-check markers
$p = $d.Paragraphs(22)
$inner = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Ability to change the active </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:rPr><w:strike/><w:lang w:val="en-US"/></w:rPr><w:t>TileMap</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' + `
    '</w:p>'
$p.Range.InsertXML((Wrap-Xml $inner))

# --- Paragraph 21: "Level objects (Tilemaps, Particles, etc)" -> split with spell-check markers
$p = $d.Paragraphs(21)
$inner = '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr>' + `
    '<w:r><w:t>Level objects (</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>Tilemaps</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve">, Particles, </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>etc</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>)</w:t></w:r>' + `
    '</w:p>'
$p.Range.InsertXML((Wrap-Xml $inner))

# --- Paragraph 1: "SuperTux Classic" -> split with spell-check markers
$p = $d.Paragraphs(1)
$inner = '<w:p><w:pPr><w:pStyle w:val="Title"/><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>SuperTux</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> Classic</w:t></w:r>' + `
    '</w:p>'
$p.Range.InsertXML((Wrap-Xml $inner))
